$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '61.667.16'
$ws.Range('E2').Value = '  -0.34%  '
$ws.Range('D3').Value = '3.368.81'
$ws.Range('E3').Value = '  -2.62%  '
$ws.Range('E4').Value = '  +0.28%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '402.47'
$ws.Range('E5').Value = '  -4.04%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '128.94'
$ws.Range('E6').Value = '  +7.51%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.592'
$ws.Range('E7').Value = '  +3.18%  '
$ws.Range('E8').Value = '  +0.11%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.662'
$ws.Range('E9').Value = '  +4.49%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.119'
$ws.Range('E10').Value = '  +4.24%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '41.30'
$ws.Range('E11').Value = '  +0.14%  '
$ws.Range('E12').Value = '  -1.05%  '
$ws.Range('D13').Value = '3.928.06'
$ws.Range('E13').Value = '  -0.89%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '8.32'
$ws.Range('E14').Value = '  -1.15%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '19.38'
$ws.Range('E15').Value = '  -1.44%  '
$ws.Range('D16').Value = '3.381.00'
$ws.Range('E16').Value = '  -1.98%  '
$ws.Range('D17').Value = '61.580.51'
$ws.Range('E17').Value = '  -0.13%  '
$ws.Range('B18').Value = 'Polygon'
$ws.Range('C18').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '1.01'
$ws.Range('E18').Value = '  -1.52%  '
$ws.Range('B19').Value = 'Uniswap'
$ws.Range('C19').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '11.16'
$ws.Range('E19').Value = '  +0.97%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '0.0000127'
$ws.Range('E20').Value = '  +9.66%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '3.22'
$ws.Range('E21').Value = '  -3.70%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '82.03'
$ws.Range('E22').Value = '  +8.97%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '12.70'
$ws.Range('E23').Value = '  +0.46%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '303.29'
$ws.Range('E24').Value = '  +2.27%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '3.11'
$ws.Range('E25').Value = '  -1.28%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '4.77'
$ws.Range('E26').Value = '  +12.29%  '
$ws.Range('B27').Value = 'Filecoin'
$ws.Range('C27').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '8.23'
$ws.Range('E27').Value = '  +6.62%  '
$ws.Range('B28').Value = 'EthereumClassic'
$ws.Range('C28').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '29.18'
$ws.Range('E28').Value = '  -7.06%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '7.49'
$ws.Range('E29').Value = '  -6.63%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '0.172'
$ws.Range('E30').Value = '  -0.44%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '0.115'
$ws.Range('E31').Value = '  +2.20%  '
$ws.Range('E32').Value = '  -0.03%  '
$ws.Range('E33').Value = '  -1.85%  '
$ws.Range('B34').Value = 'InjectiveProtocol'
$ws.Range('C34').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '41.02'
$ws.Range('E34').Value = '  -8.26%  '
$ws.Range('B35').Value = 'Toncoin'
$ws.Range('C35').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '2.48'
$ws.Range('E35').Value = '  -0.94%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.0479'
$ws.Range('E36').Value = '  -0.36%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '52.01'
$ws.Range('E37').Value = '  -0.63%  '
$ws.Range('E38').Value = '  +0.39%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '3.38'
$ws.Range('E39').Value = '  -3.42%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '2.94'
$ws.Range('E40').Value = '  -3.88%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '137.79'
$ws.Range('E41').Value = '  +2.47%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '1.96'
$ws.Range('E42').Value = '  -0.66%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.123'
$ws.Range('E43').Value = '  +1.55%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.292'
$ws.Range('E44').Value = '  +2.23%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '3.89'
$ws.Range('E45').Value = '  -1.71%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '16.60'
$ws.Range('E46').Value = '  -5.13%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '2.21'
$ws.Range('E47').Value = '  -0.34%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '21.12'
$ws.Range('E48').Value = '  -2.47%  '
$ws.Range('D49').Value = '3.702.08'
$ws.Range('E49').Value = '  -1.34%  '
$ws.Range('D50').Value = '2.105.41'
$ws.Range('E50').Value = '  -4.38%  '
$ws.Range('E51').Value = '  -4.64%  '
